# 4.1.1.1a.xlsx — add a new "2021" column (T) mirroring the existing
# per-year columns (D:S). Formatting for column T is copied from the
# adjacent 2020 column (S) so fonts/borders/number formats/wrap stay in
# sync, then data cells are switched to right-aligned (matching how the
# source workbook renders every other yearly data column).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- formatting -----------------------------------------------------
# Rows 3-5 and 23 are header/spacer rows whose S-column formatting is
# copied verbatim (no horizontal-alignment override).
$ws.Range("S3:S5").Copy()
$ws.Range("T3:T5").PasteSpecial(-4122)

$ws.Range("S23").Copy()
$ws.Range("T23").PasteSpecial(-4122)

# Data rows: copy format from column S, then force right alignment
# (mirrors columns D:S, whose data cells are right aligned).
$ws.Range("S6:S22").Copy()
$ws.Range("T6:T22").PasteSpecial(-4122)
$ws.Range("T6:T22").HorizontalAlignment = -4152

$ws.Range("S24:S40").Copy()
$ws.Range("T24:T40").PasteSpecial(-4122)
$ws.Range("T24:T40").HorizontalAlignment = -4152

$excel.CutCopyMode = 0

# --- values -----------------------------------------------------------
$ws.Range("T4").Value = 2021

$ws.Range("T6").Value = 1466
$ws.Range("T7").Value = ""
$ws.Range("T8").Value = 76
$ws.Range("T9").Value = 15
$ws.Range("T10").Value = 1
$ws.Range("T11").Value = 188
$ws.Range("T12").Value = 22
$ws.Range("T13").Value = 15
$ws.Range("T14").Value = "-"
$ws.Range("T15").Value = "-"
$ws.Range("T16").Value = 112
$ws.Range("T17").Value = "-"
$ws.Range("T18").Value = 6
$ws.Range("T19").Value = "-"
$ws.Range("T20").Value = 29
$ws.Range("T21").Value = 1002
$ws.Range("T22").Value = "-"

$ws.Range("T24").Value = 1029
$ws.Range("T25").Value = ""
$ws.Range("T26").Value = 51
$ws.Range("T27").Value = 4
$ws.Range("T28").Value = "-"
$ws.Range("T29").Value = 127
$ws.Range("T30").Value = 14
$ws.Range("T31").Value = 12
$ws.Range("T32").Value = "-"
$ws.Range("T33").Value = "-"
$ws.Range("T34").Value = 70
$ws.Range("T35").Value = "-"
$ws.Range("T36").Value = 3
$ws.Range("T37").Value = "-"
$ws.Range("T38").Value = 16
$ws.Range("T39").Value = 732
$ws.Range("T40").Value = "-"

# --- selection ----------------------------------------------------------
$ws.Range("T3").Select()
